$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new "status_label" column before column B.
# This shifts the existing B:I columns (NCTId..results) one to the
# right, to C:J, without touching their values/types - which keeps
# already-correct cell typing (numbers-as-text, multi-line text, etc.)
# intact for rows that do not otherwise change.
$ws.Columns.Item(2).Insert()

# Header row
$ws.Range("A1").Value = 'statut'
$ws.Range("B1").Value = 'status_label'
$ws.Range("C1").Value = 'NCTId'
$ws.Range("D1").Value = 'eudraCT'
$ws.Range("E1").Value = 'completion_year'
$ws.Range("F1").Value = 'clinical_trial_title'
$ws.Range("G1").Value = 'acronym'
$ws.Range("H1").Value = 'results_1y'
$ws.Range("I1").Value = 'results_3y'
$ws.Range("J1").Value = 'results'

# Fix up the status emoji (column A) for rows whose status changed
# as part of this edit (their data moved to/from another row).
$ws.Range("A9").Value = '🟥'
$ws.Range("A10").Value = '🟩'
$ws.Range("A12").Value = '🟥'
$ws.Range("A13").Value = '🟥'

# New "status_label" column: derive from the status emoji already in column A
# (rouge = red square, vert = green square).
$ws.Range("B2").Value = 'rouge'
$ws.Range("B3").Value = 'rouge'
$ws.Range("B4").Value = 'rouge'
$ws.Range("B5").Value = 'rouge'
$ws.Range("B6").Value = 'rouge'
$ws.Range("B7").Value = 'rouge'
$ws.Range("B8").Value = 'rouge'
$ws.Range("B9").Value = 'rouge'
$ws.Range("B10").Value = 'vert'
$ws.Range("B11").Value = 'rouge'
$ws.Range("B12").Value = 'rouge'
$ws.Range("B13").Value = 'rouge'
$ws.Range("B14").Value = 'rouge'
$ws.Range("B15").Value = 'rouge'
$ws.Range("B16").Value = 'rouge'

# Rows whose data moved between rows: rewrite the full row (C:J) to its
# new content.
# Row 9
$ws.Range("C9").Value = 'NCT03934073'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '2021'
$ws.Range("F9").Value = 'Evaluation of the Benefit of the Training of the Manual Dexterity Post-stroke: Effect on the Function and Spontaneous Use of the Hand and the Cerebral Plasticity'
$ws.Range("G9").Value = 'DEXTRAIN'
$ws.Range("H9").Value = $false
$ws.Range("I9").Value = $false
$ws.Range("J9").Value = $false

# Row 10
$ws.Range("C10").Value = 'NCT04350580'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '2021'
$ws.Range("F10").Value = 'Value of Early Treatment With Polyvalent Immunoglobulin in the Management of Acute Respiratory Distress Syndrome Associated With SARS-CoV-2 Infections'
$ws.Range("G10").Value = 'ICAR'
$ws.Range("H10").Value = $true
$ws.Range("I10").Value = $true
$ws.Range("J10").Value = $true

# Row 12
$ws.Range("C12").Value = 'NCT02476435'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '2022'
$ws.Range("F12").Value = 'Depersonalization Disorder: Therapeutic Effect of Neuronavigated Repetitive Transcranial Stimulation of Right Angular Gyrus.'
$ws.Range("G12").Value = 'PERSONA'
$ws.Range("H12").Value = $false
$ws.Range("I12").Value = $false
$ws.Range("J12").Value = $false

# Row 13
$ws.Range("C13").Value = 'NCT02235012'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '2022'
$ws.Range("F13").Value = 'Cognitive Biases in Decision Making in a Pharmacological Model of Psychosis : a Study in Healthy Humans Recieving Low Dose Anesthetic, Ketamine Versus Placebo'
$ws.Range("G13").Value = 'KETABI'
$ws.Range("H13").Value = $false
$ws.Range("I13").Value = $false
$ws.Range("J13").Value = $false

